$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: target cluster changes from "ECs" to "MuSCs"; dependent specificity/weight
# columns are recomputed for the new (Ccl21b/Ccr10, FAPs -> MuSCs) pairing.
$ws.Range("D2").Value = "MuSCs"
$ws.Range("I2").Value = 0.4877525841056716
$ws.Range("J2").Value = 0.588184597482006
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.1452775
$ws.Range("N2").Value = 0.290555
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.05482884229416667
$ws.Range("R2").Value = 0.328973053765
$ws.Range("S2").Value = 0.4877525841056716
$ws.Range("T2").Value = 0.588184597482006

# Row 3: sending cluster changes from "FAPs" to "MuSCs", receptor-expressing cell
# count drops from 3 to 2, and the dependent columns are recomputed accordingly.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.396361
$ws.Range("H3").Value = 0.7927219999999999
$ws.Range("I3").Value = 0.5122474158943284
$ws.Range("J3").Value = 0.411815402517994
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.1452775
$ws.Range("N3").Value = 0.290555
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.0575823351775
$ws.Range("R3").Value = 0.23032934071
$ws.Range("S3").Value = 0.5122474158943284
$ws.Range("T3").Value = 0.411815402517994

# Rows 4 and 5 (the old ECs-target and MuSCs-target pairs sourced from MuSCs)
# are no longer part of the updated TPM output, so remove them. This also
# drops the now-unused "ECs" shared string and shrinks the dimension to T3.
$ws.Rows("4:5").Delete() | Out-Null
